$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new cell A2 value (will go into shared strings as "TC3")
$ws.Range("A2").Value = "TC3"

# Update the selection to A3 (matches sheetView selection activeCell="A3" sqref="A3")
$ws.Range("A3").Select()
